$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells stay text-formatted so numeric-looking strings
# (e.g. "645.00", "11.70", thousands like "95.929.68") keep their
# exact literal formatting instead of being coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.929.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.661.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +10.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.37"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "645.00"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.04%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.402"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.659.11"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +9.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.95"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.40%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.358.57"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +10.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.764.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.17%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.58%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.46"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.76%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.667.47"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +10.26%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.48"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +23.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.72"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +7.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "519.63"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.44"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.481"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +9.37%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.70%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "93.61"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.16"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +20.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.70"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.99"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +16.38%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.22%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.581"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "561.55"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.53%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.969"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +11.23%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +7.02%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.74%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.94"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +50.80%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.07"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.33"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.50"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.77%  "
